$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja9")
$v = $ws.Cells.Item(51, 2).Value
$ws.Cells.Item(20, 15).Value = "TYPE:" + $v.GetType().FullName + ":" + $v.ToString()
